# Adds two new slides (37 "Blood Test Results - Data Decoupling" and
# 38 "AWS S3") to the end of the deck, right after the existing last
# slide ("Blood Test Results - Segmentation (4)"). Both new slides reuse
# that slide's layout/placeholder geometry, so we build them by
# duplicating it and then rewriting the title + body text.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 37: "Blood Test Results - Data Decoupling"
# ---------------------------------------------------------------------
$sourceSlide = $p.Slides.Item($p.Slides.Count)
$dup37 = $sourceSlide.Duplicate()
$s37 = $dup37.Item(1)

$title37 = $s37.Shapes.Item(2)
$title37.TextFrame.TextRange.Runs(2).Text = "Results " + [char]0x2013 + " Data Decoupling"

$body37 = $s37.Shapes.Item(1)
$tr37 = $body37.TextFrame.TextRange
$tr37.Text = "Store patient data in cloud storage`r" + `
             "blood-results.json`r" + `
             "AWS S3 - Simple Storage Service`r" + `
             "Data stored in buckets`r" + `
             "X"

# Paragraph 3 ("AWS S3 - Simple Storage Service") moves up to the top level.
$tr37.Paragraphs(3).IndentLevel = 1

# Paragraph 4 ("Data stored in buckets") becomes a level-1 sub-bullet at
# the smaller sub-bullet font size.
$tr37.Paragraphs(4).IndentLevel = 2
$tr37.Paragraphs(4).Font.Size = 30.5

# Paragraph 5 is a trailing, empty level-1 sub-bullet line. Give it the
# right level/size, then remove the temporary placeholder character.
$p5 = $tr37.Paragraphs(5)
$p5.IndentLevel = 2
$p5.Font.Size = 30.5
$tr37.Characters($p5.Start, 1).Delete()

# ---------------------------------------------------------------------
# Slide 38: "AWS S3"
# ---------------------------------------------------------------------
$dup38 = $s37.Duplicate()
$s38 = $dup38.Item(1)

$title38 = $s38.Shapes.Item(2)
$title38.TextFrame.TextRange.Text = "AWS S3"

$body38 = $s38.Shapes.Item(1)
$tr38 = $body38.TextFrame.TextRange
$tr38.Text = "Log into the AWS console`r" + `
             "Create an S3 bucket`r" + `
             "Bucket names must be globally unique!`r" + `
             ("Upload file from data directory " + [char]0x201C + "blood-results.json" + [char]0x201D) + "`r" + `
             "Change IAM role to allow Lambda function to read S3 bucket (or create a new role)"

# Append the two remaining sub-bullets (IAM policy names) after the last
# paragraph created above.
$lastPara = $tr38.Paragraphs($tr38.Paragraphs().Count)
$lastPara.InsertAfter("`rlambda_basic_execution`rAmazonS3ReadOnlyAccess")

# Paragraph 2 ("Create an S3 bucket") stays at the top level (size 32).
$tr38.Paragraphs(2).IndentLevel = 1
$tr38.Paragraphs(2).Font.Size = 32

# Paragraph 3 ("Bucket names must be globally unique!") is a smaller
# level-1 sub-bullet.
$tr38.Paragraphs(3).IndentLevel = 2
$tr38.Paragraphs(3).Font.Size = 29

# Paragraph 4 (upload instructions) is a level-1 sub-bullet.
$tr38.Paragraphs(4).IndentLevel = 2
$tr38.Paragraphs(4).Font.Size = 30.5

# Paragraph 5 (Change IAM role...) stays at the top level (size 32).
$tr38.Paragraphs(5).IndentLevel = 1
$tr38.Paragraphs(5).Font.Size = 32

# Paragraphs 6 & 7 (IAM policy names) are level-1 sub-bullets.
$tr38.Paragraphs(6).IndentLevel = 2
$tr38.Paragraphs(6).Font.Size = 30.5
$tr38.Paragraphs(7).IndentLevel = 2
$tr38.Paragraphs(7).Font.Size = 30.5

Write-Output ("Slides now: " + $p.Slides.Count)
